# Weekly update: insert a new price record as the first data row (row 50)
# for "Agrícola del Norte S.A. de Arica - Mandarina", pushing the existing
# rows 50-109 down by one (to rows 51-110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 50; Excel shifts rows 50..109 down to 51..110.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with this week's record.
$ws.Cells.Item(50, 1).Value  = 1
$ws.Cells.Item(50, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value  = 44763
$ws.Cells.Item(50, 5).Value  = 15
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100102
$ws.Cells.Item(50, 8).Value  = "Cítricos"
$ws.Cells.Item(50, 9).Value  = 100102004
$ws.Cells.Item(50, 10).Value = "Mandarina"
$ws.Cells.Item(50, 11).Value = "Murcott"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 250
$ws.Cells.Item(50, 14).Value = 12000
$ws.Cells.Item(50, 15).Value = 13000
$ws.Cells.Item(50, 16).Value = 12500
$ws.Cells.Item(50, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(50, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(50, 19).Value = 625
$ws.Cells.Item(50, 20).Value = 20
